# Update the worksheet date and the 25 division-problem answers in the
# 5x5 results table. Cells are addressed by (row, column) in the Word
# table object model so the edit is unambiguous even where some of the
# new values happen to collide textually with other (unrelated) old
# values elsewhere in the document.

$d = $word.ActiveDocument

# --- Title / date paragraph -------------------------------------------------
$d.Paragraphs.Item(1).Range.Text = "2024-06-29 Saturday"

# --- Table of answers --------------------------------------------------------
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "24÷8=3, 0" },
    @{ Row = 1;  Col = 2; Text = "94÷8=11, 6" },
    @{ Row = 1;  Col = 3; Text = "17÷9=1, 8" },
    @{ Row = 1;  Col = 4; Text = "48÷2=24, 0" },
    @{ Row = 1;  Col = 5; Text = "78÷7=11, 1" },

    @{ Row = 5;  Col = 1; Text = "64÷3=21, 1" },
    @{ Row = 5;  Col = 2; Text = "34÷6=5, 4" },
    @{ Row = 5;  Col = 3; Text = "63÷9=7, 0" },
    @{ Row = 5;  Col = 4; Text = "13÷7=1, 6" },
    @{ Row = 5;  Col = 5; Text = "38÷7=5, 3" },

    @{ Row = 9;  Col = 1; Text = "72÷2=36, 0" },
    @{ Row = 9;  Col = 2; Text = "15÷2=7, 1" },
    @{ Row = 9;  Col = 3; Text = "21÷9=2, 3" },
    @{ Row = 9;  Col = 4; Text = "18÷5=3, 3" },
    @{ Row = 9;  Col = 5; Text = "35÷4=8, 3" },

    @{ Row = 13; Col = 1; Text = "78÷9=8, 6" },
    @{ Row = 13; Col = 2; Text = "29÷6=4, 5" },
    @{ Row = 13; Col = 3; Text = "20÷9=2, 2" },
    @{ Row = 13; Col = 4; Text = "39÷3=13, 0" },
    @{ Row = 13; Col = 5; Text = "65÷5=13, 0" },

    @{ Row = 17; Col = 1; Text = "24÷4=6, 0" },
    @{ Row = 17; Col = 2; Text = "37÷2=18, 1" },
    @{ Row = 17; Col = 3; Text = "72÷7=10, 2" },
    @{ Row = 17; Col = 4; Text = "33÷5=6, 3" },
    @{ Row = 17; Col = 5; Text = "21÷8=2, 5" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
